$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.061.73"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "1.825.96"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.31"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07230"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8444"
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("D12").Value = "1.830.08"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.290"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07032"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.71"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008743"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.88"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").Value = "27.103.55"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.128"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.81"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").Value = "2.052.38"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.980"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.44"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.262"
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.237"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.69"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08742"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7367"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.431"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.095"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05235"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.317"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.871"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1687"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5102"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.567"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.981"
$ws.Range("E45").Value = "  +7.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.52"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4728"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.51"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06322"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.652"
$ws.Range("E51").Value = "  -1.64%  "
